$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.645.71"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "2.078.88"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "232.46"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "58.29"
$ws.Range("E7").Value = "  -1.34%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.388"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").Value = "0.0778"
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").Value = "2.383.36"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "14.82"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "21.38"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").Value = "0.769"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "2.075.89"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "37.627.73"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").Value = "6.15"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "70.03"
$ws.Range("E20").Value = "  -2.05%  "
$ws.Range("D21").Value = "0.0₃0828"
$ws.Range("E21").Value = "  -2.19%  "
$ws.Range("D22").Value = "227.83"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("D26").Value = "9.89"
$ws.Range("E26").Value = "  +7.18%  "
$ws.Range("D27").Value = "169.94"
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("E28").Value = "  -3.63%  "
$ws.Range("D29").Value = "19.31"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("E30").Value = "  -3.94%  "
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("D32").Value = "4.58"
$ws.Range("E32").Value = "  -3.24%  "
$ws.Range("D33").Value = "0.0628"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("D34").Value = "4.64"
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("D35").Value = "2.54"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").Value = "3.32"
$ws.Range("E37").Value = "  -3.22%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "5.33"
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("E40").Value = "  +4.18%  "
$ws.Range("D41").Value = "99.13"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  -1.92%  "
$ws.Range("E43").Value = "  +4.63%  "
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").Value = "1.482.49"
$ws.Range("E45").Value = "  +2.46%  "
$ws.Range("D46").Value = "16.82"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("E48").Value = "  -3.35%  "
$ws.Range("D49").Value = "7.26"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("D50").Value = "2.98"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").Value = "2.267.90"
$ws.Range("E51").Value = "  -0.35%  "
